# Add a new data row (row 2) under the existing header row, mirroring a
# fresh "test data" entry. All values are entered as plain text (matching
# the workbook's existing convention of storing everything, including
# dates and numbers, as text) rather than being auto-converted by Excel
# into a date serial / numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to Text format first so Excel doesn't
# auto-convert "09/08/2023" into a date or "5000.00" into a number.
$row2 = $ws.Range("A2:G2")
$row2.NumberFormat = "@"

$ws.Range("A2").Value = "09/08/2023"
$ws.Range("B2").Value = "5000.00"
$ws.Range("C2").Value = "5000.00"
$ws.Range("D2").Value = "6000.00"
$ws.Range("E2").Value = "6000.00"
$ws.Range("F2").Value = "1000.00"
$ws.Range("G2").Value = "120.00"

# Drop back to the default (unstyled) cell style for this row, same as
# the header row's formatting was never intended to apply here.
$row2.Style = "Normal"
